$wb = $excel.ActiveWorkbook

# Rename the existing sheet and add the new "Solubility" sheet right after it.
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Inventory"
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Solubility"

# Header row - written in the same order the original author typed them
# (A, B, D, C, E) so the shared-string table ends up in the matching order.
$ws2.Range("A1").Value = "Component 1 Abbreviation"
$ws2.Range("B1").Value = "Solvent Abbreviation"
$ws2.Range("D1").Value = "Temperature °C"
$ws2.Range("C1").Value = "Component 1 Maximum g/L"
$ws2.Range("E1").Value = "Source"

# Row 2 - Sodium Chloride in Water
$ws2.Range("A2").Value = "NaCl"
$ws2.Range("B2").Value = "Water"
$ws2.Range("C2").Value = 357
$ws2.Range("D2").Value = 25
$ws2.Hyperlinks.Add($ws2.Range("E2"), "https://www.sigmaaldrich.com/content/dam/sigma-aldrich/docs/Sigma-Aldrich/Product_Information_Sheet/s7653pis.pdf")

# Row 3 - Butenedioic acid in Water
$ws2.Range("A3").Value = "Butenedioicacid"
$ws2.Range("B3").Value = "Water"
$ws2.Range("C3").Value = 441
$ws2.Range("D3").Value = 25
$rng = $ws2.Range("E3")
$rng.Value = "Yalkowsky, S.H., He, Yan, Jain, P. Handbook of Aqueous Solubility Data Second Edition. CRC Press, Boca Raton, FL 2010, p. 86"
$f = $rng.Font
$f.Name = "Segoe UI"
$f.Size = 9.6
$f.Italic = $true
$f.Color = 7037275
$rng.HorizontalAlignment = -4131
$rng.VerticalAlignment = -4108
$rng.WrapText = $true
$rng.IndentLevel = 1

# Row 4 - Choline chloride in Water
$ws2.Range("A4").Value = "chcl"
$ws2.Range("B4").Value = "Water"
$ws2.Range("C4").Value = 650

# Row 5 - Gold (III) Chloride Trihydrate in Water
$ws2.Range("A5").Value = "haucl4"
$ws2.Range("B5").Value = "Water"
$ws2.Range("C5").Value = 150
$ws2.Hyperlinks.Add($ws2.Range("E5"), "https://www.emdmillipore.com/US/en/product/TetrachloroauricIII-acid-trihydrate-990-0,MDA_CHEM-101582?ReferrerURL=https%3A%2F%2Fwww.google.com%2F&bd=1")

# Column widths for the Solubility sheet
$ws2.Columns.Item(1).ColumnWidth = 23.08984375
$ws2.Columns.Item(2).ColumnWidth = 23.08984375
$ws2.Columns.Item(3).ColumnWidth = 23.453125
$ws2.Columns.Item(4).ColumnWidth = 23.08984375
$ws2.Columns.Item(5).ColumnWidth = 104.54296875

# View state: Inventory scrolled with B2:B6 selected; Solubility tab active with C6 selected.
$ws1.Application.ActiveWindow.ScrollRow = 2
$ws1.Range("B2:B6").Select()
$ws2.Range("C6").Select()
$ws2.Activate()

Write-Host "done"
